$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Bring "Repayment schedule" to the front (sets activeTab in workbook.xml and
# tabSelected on this sheet's view, clearing it from whatever was active
# before - the "Summary" sheet).
$ws.Activate()

# Insert a new (blank) column before column N, pushing the old N/O/P
# ("Late" / blank / "Outstanding") columns one to the right.
$ws.Columns("N").Insert()

# The newly inserted column inherits the width of the column to its left
# (M - "In Advance"), matching Excel's normal "insert column" behaviour.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Move the selection to R7, as recorded in the saved view.
$ws.Range("R7").Select()
